$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handback transform failed for the 7df61d27-... item: update the
# "Status" text everywhere it is shown (Overview summary row + each
# language sheet's Status column) for that file.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file-name mismatch in the Error Detail
# column (P) for the 7df61d27-... row on each language sheet.
$wsZhCn.Range("P3").Value = "Handback file name: xkne3olh.hkn is different with handoff file name: 7df61d27-f421-48f2-8d66-93250073fc34.d3a8773048f660ed4d3bc5c37f656bac80270511.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: xkne3olh.hkn is different with handoff file name: 7df61d27-f421-48f2-8d66-93250073fc34.d3a8773048f660ed4d3bc5c37f656bac80270511.de-de."

# Widen the Error Detail column so the new message is readable.
$wsZhCn.Range("P1").ColumnWidth = 39.15
$wsDeDe.Range("P1").ColumnWidth = 39.15
